$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$newValues = @{
    2  = "hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
    3  = "Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
    4  = "Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
    5  = "Veronia.rafat@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, hananragab@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
    6  = "Veronia.rafat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
    13 = "wessam.atef@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg"
    24 = "Remon.Matta@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, marina_atef@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"
    25 = "Remon.Matta@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, marina_atef@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
    30 = "Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
    31 = "Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
    32 = "Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
    33 = "Veronia.rafat@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, hananragab@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
    34 = "Veronia.rafat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
    41 = "wessam.atef@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg"
    52 = "Remon.Matta@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, marina_atef@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"
    53 = "Remon.Matta@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, marina_atef@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
